# model_7_date_2018-05-14_21-44-21__epochs_100_step_3.xlsx
# - rename the summary-row label from "min" to "min/max"
# - add MAX/MIN roll-up formulas for columns B (acc) and C (loss), mirroring
#   the existing D/E roll-ups
# - scroll the view down and move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. "min" -> "min/max" on the summary row label (A53, shared string)
$ws.Range("A53").Value = "min/max"

# 2. New roll-up formulas on row 53: MAX for acc (B), MIN for loss (C) -
#    same pattern already used for D53 (MAX) / E53 (MIN)
$ws.Range("B53").Formula = "=MAX(B1:B51)"
$ws.Range("C53").Formula = "=MIN(C1:C51)"

# 3. View state: scroll so row 40 is at the top, and move the selection to C54
$win = $excel.ActiveWindow
$win.ScrollRow = 40
$win.ScrollColumn = 1
$ws.Range("C54").Select()
